$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: insert "revenue" before "match", drop "email"
$ws.Range("E1").Value = "revenue"
$ws.Range("F1").Value = "match"

# Replace row 2 with Dixon Technologies data
$ws.Range("A2").Value = "Dixon Technologies"
$ws.Range("B2").Value = "https://dixonindia.com"
$ws.Range("C2").Value = "Dixon Technologies is an electronics manufacturing services (EMS) company"
$ws.Range("D2").Value = "LED TVs, mobile phones, home appliances"
$ws.Range("E2").Value = "500-1000 cr"
$ws.Range("F2").Value = "Replicant Systems' vision AI can help Dixon Technologies automate inspection of LED TVs and mobile phones, reducing defects and increasing efficiency."

# Add row 3 with Optiemus Infracom data
$ws.Range("A3").Value = "Optiemus Infracom"
$ws.Range("B3").Value = "https://optiemus.com"
$ws.Range("C3").Value = "Optiemus Infracom is a telecom and electronics manufacturing company"
$ws.Range("D3").Value = "Telecom equipment, electronic devices, IoT solutions"
$ws.Range("E3").Value = "500-1000 cr"
$ws.Range("F3").Value = "Replicant Systems' industrial automation solutions can assist Optiemus Infracom in streamlining telecom equipment manufacturing, improving product reliability and reducing costs."
